# "Added Member Type and Salutation in Member Form"
# Rows 7, 10 and 12 move from "In-Progress" (Add-On/Enhancement still pending)
# to "Completed" - matching the look (fill/number formats) already used by
# the other completed rows (e.g. row 8), plus a completion date and full
# (100%) completion percentage. Row 12's Man Hours is corrected from 6 to 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of an already-"Completed" row (row 8) onto the rows
# that are being marked completed, so the cell styles (fill color, date /
# percent number formats, borders, comma formatting on the amount column)
# match the rest of the "Completed" rows exactly.
$fmtSource = $ws.Range("A8:K8")
$fmtSource.Copy()
$ws.Range("A7:K7").PasteSpecial(-4122)
$ws.Range("A10:K10").PasteSpecial(-4122)
$ws.Range("A12:K12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 7: "Club Email, WebSite, Waiter Name should be print in Bill" -> Completed
$ws.Range("G7").Value = "Completed"
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 45168
$ws.Range("I7").NumberFormat = "d-mmm-yy"

# Row 8: "Update Membership Number" already Completed - just stamp the
# completion date that was missing.
$ws.Range("I8").Value = 45168
$ws.Range("I8").NumberFormat = "d-mmm-yy"

# Row 10: "Add Member Type in Member Creation" -> Completed
$ws.Range("G10").Value = "Completed"
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 45168
$ws.Range("I10").NumberFormat = "d-mmm-yy"

# Row 12: "Salutation in Member Creation Form" -> Completed, and its Man
# Hours estimate is corrected from 6 to 4 (the Calculated Amount formula
# recalculates automatically: 4 * 600 = 2400).
$ws.Range("F12").Value = 4
$ws.Range("G12").Value = "Completed"
$ws.Range("H12").Value = 1
$ws.Range("I12").Value = 45168
$ws.Range("I12").NumberFormat = "d-mmm-yy"

# Keep the on-disk cursor position in sync with where the edit happened.
[void]$ws.Range("C12").Select()
